$newStrings = @(
    'Sát Phá Tham hội chiếu tại Tật Ách',
    'Tử Phủ Vũ Tướng hội chiếu tại Tật Ách',
    'Cơ Nguyệt Đồng Lương hội chiếu tại Tật Ách',
    'Cự Nhật hội chiếu tại Tật Ách',
    'Kình Đà hội chiếu tại Tật Ách',
    'Xương Khúc hội chiếu tại Tật Ách',
    'Hoả Linh hội chiếu tại Tật Ách',
    'Không Kiếp hội chiếu tại Tật Ách',
    'Quang Quý hội chiếu tại Tật Ách',
    'Tả Hữu hội chiếu tại Tật Ách',
    'Song Hao hội chiếu tại Tật Ách',
    'Tang Hổ hội chiếu tại Tật Ách',
    'Khốc Hư hội chiếu tại Tật Ách',
    'Hình Riêu hội chiếu tại Tật Ách',
    'Thai Toạ hội chiếu tại Tật Ách',
    'Đào Hồng hội chiếu tại Tật Ách',
    'Ấn Phù hội chiếu tại Tật Ách',
    'Song Hao Quyền Lộc Kiếp Hoả hội chiếu tại Tật Ách',
    'Tử Phủ Vũ Tướng Xương Khúc Khôi Việt Tả Hữu Khoa Quyền Lộc Long hội chiếu tại Tật Ách',
    'Tử Phủ Vũ Tướng Tả Hữu Khoa Quyền Lộc Long Phượng hội chiếu tại Tật Ách',
    'Tử Khúc Phá Dương Đà hội chiếu tại Tật Ách',
    'Cơ Nguyệt Đồng Lương gặp Xương Khúc Tả Hữu hội chiếu tại Tật Ách',
    'Cơ Nguyệt Đồng Lương Khoa Tả Hữu Quang Quý Quan Phúc hội chiếu tại Tật Ách',
    'Sát Quyền hội chiếu tại Tật Ách',
    'Lộc Mã hội chiếu tại Tật Ách',
    'Kiếp Hư Hao Quyền hội chiếu tại Tật Ách',
    'Tuế Hổ Phù Xương Khúc hội chiếu tại Tật Ách',
    'Xương Khúc Tấu Long Phượng hội chiếu tại Tật Ách',
    'Đào Hồng Riêu Tấu Cơ Vũ hội chiếu tại Tật Ách',
    'Binh Hình Tướng Ấn hội chiếu tại Tật Ách',
    'Hổ Tấu hội chiếu tại Tật Ách',
    'Hình Riêu Y hội chiếu tại Tật Ách',
    'Mã Hỏa Linh hội chiếu tại Tật Ách',
    'Thai Tọa Hồng Đào hội chiếu tại Tật Ách',
    'Tả Hữu Không Kiếp hội chiếu tại Tật Ách',
    'Tả Hữu Binh Tướng hội chiếu tại Tật Ách',
    'Đào Quyền hội chiếu tại Tật Ách',
    'Đào Hồng Tả Cái Hữu hội chiếu tại Tật Ách',
    'Quan Phúc Quang Tấu hội chiếu tại Tật Ách',
    'Đào Hồng Xương Khúc hội chiếu tại Tật Ách',
    'Đào Tử Phủ hội chiếu tại Tật Ách',
    'Hổ Kình Sát hội chiếu tại Tật Ách',
    'Hổ Tang Kiếp hội chiếu tại Tật Ách',
    'Hổ Tang Không Kiếp hội chiếu tại Tật Ách',
    'Âm Dương Lương hội chiếu tại Tật Ách',
    'Cơ Lương Gia Hội hội chiếu tại Tật Ách',
    'Nhật Chiếu Lôi Môn hội chiếu tại Tật Ách',
    'Tả Hữu Xương Khúc hội chiếu tại Tật Ách',
    'Tham Linh Triều Viên hội chiếu tại Tật Ách',
    'Cự Hỏa Linh hội chiếu tại Tật Ách',
    'Hồng Đào Kỵ hội chiếu tại Tật Ách',
    'Cự Đồng Hình hội chiếu tại Tật Ách',
    'Lương Phá hội chiếu tại Tật Ách',
    'Khôi Việt hội chiếu tại Tật Ách',
    'Tham Vũ Hỏa hội chiếu tại Tật Ách',
    'Cự Kỵ hội chiếu tại Tật Ách',
    'Kình Đà Hỏa hội chiếu tại Tật Ách',
    'Kình Đà Không Kiếp hội chiếu tại Tật Ách',
    'Đào Hồng Hỷ hội chiếu tại Tật Ách',
    'Đào Hồng Kiếp Sát hội chiếu tại Tật Ách',
    'Đào Hồng Kỵ hội chiếu tại Tật Ách',
    'Đào Hồng Tả Phù Hữu Bật hội chiếu tại Tật Ách',
    'Đào Hồng Tả Phù Hữu Bật Khoa Quyền Lộc hội chiếu tại Tật Ách',
    'Đào Hồng Tả Phù Hữu Bật Khoa Quyền Lộc Long Trì Phượng Các hội chiếu tại Tật Ách',
    'Lưu Hà Kiếp Sát hội chiếu tại Tật Ách',
    'Phục Không Kiếp hội chiếu tại Tật Ách',
    'Xương Khúc Khôi Việt hội chiếu tại Tật Ách',
    'Khoa Quyền hội chiếu tại Tật Ách',
    'Lộc Quyền hội chiếu tại Tật Ách',
    'Tướng Binh Đào Hồng hội chiếu tại Tật Ách',
    'Lương Khốc Tuế hội chiếu tại Tật Ách',
    'Tướng Binh hội chiếu tại Tật Ách',
    'Cự Tang hội chiếu tại Tật Ách',
    'Cự Hỏa hội chiếu tại Tật Ách',
    'Phá Hình Kỵ hội chiếu tại Tật Ách',
    'Tang Trực Tuế hội chiếu tại Tật Ách',
    'Tuế Xương Khúc hội chiếu tại Tật Ách',
    'Tử Vi Tả Hữu hội chiếu tại Tật Ách',
    'Tử Vi Tang Tả Hữu hội chiếu tại Tật Ách'
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 4249
for ($i = 0; $i -lt $newStrings.Length; $i++) {
    $row = $startRow + $i
    $val = $newStrings[$i]
    $ws.Cells.Item($row, 1).Value = $val
    $ws.Cells.Item($row, 2).Value = $val
}

$lastRow = $startRow + $newStrings.Length - 1
$ws.Range("B$startRow`:B$lastRow").Select()
$ws.Application.ActiveWindow.ScrollRow = 4220
